$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "./file2.cpp(8)"
$ws.Range("B4").Value = "osx_source"
$ws.Range("A5").Value = "./file2.cpp(8)"
$ws.Range("B5").Value = "check_error"
$ws.Range("C5").Value = 2

$ws.Range("B4").Font.Name = "Arial Unicode MS"
$ws.Range("B4").Font.Size = 10
$ws.Range("B4").Font.Color = 0

$ws.Rows.Item(4).RowHeight = 17

$ws.Range("C6").Select()
